$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.490.59'
$ws.Range('E2').Value = '  +6.04%  '
$ws.Range('D3').Value = '3.549.95'
$ws.Range('E3').Value = '  +6.05%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '190.61'
$ws.Range('E5').Value = '  +9.87%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '560.61'
$ws.Range('E6').Value = '  +5.20%  '
$ws.Range('D7').Value = '3.545.56'
$ws.Range('E7').Value = '  +5.84%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.618'
$ws.Range('E8').Value = '  +3.65%  '
$ws.Range('E9').Value = '  +0.00%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.634'
$ws.Range('E10').Value = '  +3.89%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.154'
$ws.Range('E11').Value = '  +14.55%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '55.35'
$ws.Range('E12').Value = '  +3.81%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000274'
$ws.Range('E13').Value = '  +6.70%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '9.39'
$ws.Range('E14').Value = '  +1.93%  '
$ws.Range('D15').Value = '4.109.78'
$ws.Range('E15').Value = '  +6.04%  '
$ws.Range('D16').Value = '3.548.07'
$ws.Range('E16').Value = '  +6.21%  '
$ws.Range('E17').Value = '  +3.68%  '
$ws.Range('B18').Value = 'Chainlink'
$ws.Range('C18').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '18.48'
$ws.Range('E18').Value = '  +5.58%  '
$ws.Range('B19').Value = 'WrappedBTC'
$ws.Range('C19').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D19').Value = '67.477.84'
$ws.Range('E19').Value = '  +5.99%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.03'
$ws.Range('E20').Value = '  +7.12%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.00'
$ws.Range('E21').Value = '  +3.43%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '427.71'
$ws.Range('E22').Value = '  +14.71%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.10'
$ws.Range('E23').Value = '  +9.32%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '85.20'
$ws.Range('E24').Value = '  +4.21%  '
$ws.Range('E25').Value = '  -0.13%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '11.09'
$ws.Range('E26').Value = '  -2.00%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.92'
$ws.Range('E27').Value = '  +8.01%  '
$ws.Range('E28').Value = '  -0.32%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '12.36'
$ws.Range('E29').Value = '  +9.16%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '9.04'
$ws.Range('E30').Value = '  +9.02%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '30.60'
$ws.Range('E31').Value = '  +5.75%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '637.03'
$ws.Range('E32').Value = '  -1.44%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '6.68'
$ws.Range('E33').Value = '  +3.66%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '11.75'
$ws.Range('E34').Value = '  +4.84%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.112'
$ws.Range('E35').Value = '  +5.03%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '60.29'
$ws.Range('E36').Value = '  +3.45%  '
$ws.Range('B37').Value = 'PEPE'
$ws.Range('C37').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D37').Value = '0.0₃0822'
$ws.Range('E37').Value = '  +12.95%  '
$ws.Range('B38').Value = 'InjectiveProtocol'
$ws.Range('C38').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '38.53'
$ws.Range('E38').Value = '  +3.87%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.148'
$ws.Range('E39').Value = '  +19.05%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.00'
$ws.Range('E40').Value = '  -0.12%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.389'
$ws.Range('E41').Value = '  +1.96%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.46'
$ws.Range('E42').Value = '  +16.05%  '
$ws.Range('D43').Value = '3.133.10'
$ws.Range('E43').Value = '  +7.02%  '
$ws.Range('E44').Value = '  +0.08%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.66'
$ws.Range('E45').Value = '  +2.28%  '
$ws.Range('B46').Value = 'ThetaToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.89'
$ws.Range('E46').Value = '  +10.56%  '
$ws.Range('B47').Value = 'ApeXProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.38'
$ws.Range('E47').Value = '  +10.20%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0421'
$ws.Range('E48').Value = '  +5.55%  '
$ws.Range('E49').Value = '  +3.98%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.132'
$ws.Range('E50').Value = '  +5.71%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '143.63'
$ws.Range('E51').Value = '  +4.52%  '
